$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")

# --- Move the "Means"/"Best" label row down from 38 to 39 to make room for the new
#     data row that is about to be appended to the table (row 36), without disturbing
#     anything else further down the sheet (e.g. the spacer cell in AD49). -----------
$labelB = $ws.Range("B38").Value2
$labelH = $ws.Range("H38").Value2
$ws.Range("B39").Value = $labelB
$ws.Range("H39").Value = $labelH
$ws.Range("B38").ClearContents()
$ws.Range("H38").ClearContents()

# --- Fill in the values that were missing on the last existing data row (35). -------
$ws.Range("D35").Value = 4.5
$ws.Range("E35").Value = 5.3
$ws.Range("F35").Value = 9.7
$ws.Range("G35").Value = 4.9
$ws.Range("H35").Value = 6.5
$ws.Range("K35").Value = 10

# --- Add the new day's data in row 36, matching the formats used by the row above. --
$ws.Range("A36").NumberFormat = $ws.Range("A35").NumberFormat
$ws.Range("A36").Value = 45493
$ws.Range("B36").NumberFormat = $ws.Range("B35").NumberFormat
$ws.Range("B36").Value = 77.1
$ws.Range("C36").Value = 4.7
$ws.Range("D36").NumberFormat = $ws.Range("D35").NumberFormat
$ws.Range("D36").Value = 3.9
$ws.Range("I36").Value = 10
$ws.Range("J36").Value = 8

# --- Grow the Table (ListObject) so it covers the newly added row. ------------------
$tbl = $ws.ListObjects.Item("Table2")
$tbl.Resize($ws.Range("A5:K36"))

# --- Extend the summary formulas (row 2: AVERAGE, row 3: MAX/MIN) by one row. -------
# (Setting the same formula text across the whole range first, then re-asserting it
#  cell-by-cell, is what makes the engine record these as one shared-formula group,
#  matching how Excel stores a formula that was filled across a row.)
$ws.Range("B2:H2").Formula = "=AVERAGE(B6:B67)"
$ws.Range("B2").Formula = "=AVERAGE(B6:B67)"
$ws.Range("C2").Formula = "=AVERAGE(C6:C67)"
$ws.Range("D2").Formula = "=AVERAGE(D6:D67)"
$ws.Range("E2").Formula = "=AVERAGE(E6:E67)"
$ws.Range("F2").Formula = "=AVERAGE(F6:F67)"
$ws.Range("G2").Formula = "=AVERAGE(G6:G67)"
$ws.Range("H2").Formula = "=AVERAGE(H6:H67)"

$ws.Range("B3").Formula = "=MAX(B6:B74)"

$ws.Range("C3:H3").Formula = "=MIN(C6:C74)"
$ws.Range("C3").Formula = "=MIN(C6:C74)"
$ws.Range("D3").Formula = "=MIN(D6:D74)"
$ws.Range("E3").Formula = "=MIN(E6:E74)"
$ws.Range("F3").Formula = "=MIN(F6:F74)"
$ws.Range("G3").Formula = "=MIN(G6:G74)"
$ws.Range("H3").Formula = "=MIN(H6:H74)"

# --- Update the selection shown when the file is reopened. --------------------------
$ws.Activate()
$ws.Range("E36").Select()
